# Update fc2_weights data block (rows 2-11, columns A-Y) per GS Algorithm iteration
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = [double]"6.850574573036283e-05"
$ws.Range("B2").Value = [double]"0.0008434053161181509"
$ws.Range("C2").Value = [double]"8.698595621581262e-09"
$ws.Range("D2").Value = [double]"4.319526851759292e-06"
$ws.Range("E2").Value = [double]"1.770770995926287e-14"
$ws.Range("F2").Value = [double]"1.954338380016907e-08"
$ws.Range("G2").Value = [double]"2.666527507244609e-05"
$ws.Range("H2").Value = [double]"3.933255054988649e-09"
$ws.Range("I2").Value = [double]"8.556725333619397e-06"
$ws.Range("J2").Value = [double]"0.0003688785945996642"
$ws.Range("K2").Value = [double]"2.231896800708455e-08"
$ws.Range("L2").Value = [double]"1.447957043154702e-08"
$ws.Range("M2").Value = [double]"0.008960417471826077"
$ws.Range("N2").Value = [double]"4.422629009366639e-17"
$ws.Range("O2").Value = [double]"6.479158764705062e-05"
$ws.Range("P2").Value = [double]"2.048904867990586e-08"
$ws.Range("Q2").Value = [double]"2.548053679163331e-08"
$ws.Range("R2").Value = [double]"1.482119937890047e-08"
$ws.Range("S2").Value = [double]"1.880398667708505e-05"
$ws.Range("T2").Value = [double]"6.298531843640376e-06"
$ws.Range("U2").Value = [double]"0.002949760761111975"
$ws.Range("V2").Value = [double]"2.668303068276145e-06"
$ws.Range("W2").Value = [double]"2.320366121466577e-08"
$ws.Range("X2").Value = [double]"0.0004497090412769467"
$ws.Range("Y2").Value = [double]"0.002099902601912618"

$ws.Range("A3").Value = [double]"0.0002923326683230698"
$ws.Range("B3").Value = [double]"0.001285611884668469"
$ws.Range("C3").Value = [double]"2.299433532471085e-09"
$ws.Range("D3").Value = [double]"1.917597414191619e-09"
$ws.Range("E3").Value = [double]"2.475292157259e-17"
$ws.Range("F3").Value = [double]"6.374432537370467e-09"
$ws.Range("G3").Value = [double]"4.16068323829677e-05"
$ws.Range("H3").Value = [double]"2.234758953912665e-17"
$ws.Range("I3").Value = [double]"0.0001317159330938011"
$ws.Range("J3").Value = [double]"3.191697760485113e-05"
$ws.Range("K3").Value = [double]"0.0005330232670530677"
$ws.Range("L3").Value = [double]"4.276933029956353e-09"
$ws.Range("M3").Value = [double]"0.00777438236400485"
$ws.Range("N3").Value = [double]"2.147962065056863e-17"
$ws.Range("O3").Value = [double]"7.201521384558873e-06"
$ws.Range("P3").Value = [double]"5.514883797026748e-13"
$ws.Range("Q3").Value = [double]"2.212343641616599e-09"
$ws.Range("R3").Value = [double]"2.145699389650632e-13"
$ws.Range("S3").Value = [double]"0.0002137597184628248"
$ws.Range("T3").Value = [double]"1.551671630295459e-05"
$ws.Range("U3").Value = [double]"2.921883424278349e-05"
$ws.Range("V3").Value = [double]"3.211889998055995e-05"
$ws.Range("W3").Value = [double]"1.66432667647598e-09"
$ws.Range("X3").Value = [double]"0.0001709639182081446"
$ws.Range("Y3").Value = [double]"0.002281186869367957"

$ws.Range("A4").Value = [double]"4.164793426753022e-05"
$ws.Range("B4").Value = [double]"0.000873119046445936"
$ws.Range("C4").Value = [double]"9.952799695511771e-10"
$ws.Range("D4").Value = [double]"2.023156639552326e-06"
$ws.Range("E4").Value = [double]"5.997534913149138e-08"
$ws.Range("F4").Value = [double]"4.433336653164588e-05"
$ws.Range("G4").Value = [double]"0.0004010562843177468"
$ws.Range("H4").Value = [double]"6.594645321200422e-18"
$ws.Range("I4").Value = [double]"0.002083619823679328"
$ws.Range("J4").Value = [double]"4.361839091870934e-05"
$ws.Range("K4").Value = [double]"4.966978041642066e-18"
$ws.Range("L4").Value = [double]"1.273284169656108e-06"
$ws.Range("M4").Value = [double]"0.01062809582799673"
$ws.Range("N4").Value = [double]"7.687747129896385e-18"
$ws.Range("O4").Value = [double]"0.000277571874903515"
$ws.Range("P4").Value = [double]"8.418908947760428e-08"
$ws.Range("Q4").Value = [double]"1.199122380057815e-05"
$ws.Range("R4").Value = [double]"8.10360045733205e-09"
$ws.Range("S4").Value = [double]"0.0001946452393895015"
$ws.Range("T4").Value = [double]"0.0002477207162883133"
$ws.Range("U4").Value = [double]"0.0001232908689416945"
$ws.Range("V4").Value = [double]"4.831540536542889e-07"
$ws.Range("W4").Value = [double]"1.573590964198956e-07"
$ws.Range("X4").Value = [double]"0.0006239789072424173"
$ws.Range("Y4").Value = [double]"0.001610737759619951"

$ws.Range("A5").Value = [double]"0.0004995528142899275"
$ws.Range("B5").Value = [double]"0.0007081272779032588"
$ws.Range("C5").Value = [double]"0.001123439986258745"
$ws.Range("D5").Value = [double]"9.34303159283445e-08"
$ws.Range("E5").Value = [double]"0.0001596799120306969"
$ws.Range("F5").Value = [double]"1.792307898540457e-06"
$ws.Range("G5").Value = [double]"0.0003022235177922994"
$ws.Range("H5").Value = [double]"3.338568532740283e-08"
$ws.Range("I5").Value = [double]"2.622188367240597e-05"
$ws.Range("J5").Value = [double]"1.333777709078277e-05"
$ws.Range("K5").Value = [double]"1.415430270606236e-17"
$ws.Range("L5").Value = [double]"5.402661145126331e-07"
$ws.Range("M5").Value = [double]"0.009919251315295696"
$ws.Range("N5").Value = [double]"4.544984903830993e-17"
$ws.Range("O5").Value = [double]"9.185194357996807e-05"
$ws.Range("P5").Value = [double]"2.218024377498296e-07"
$ws.Range("Q5").Value = [double]"3.165545786032453e-05"
$ws.Range("R5").Value = [double]"1.02748308563605e-05"
$ws.Range("S5").Value = [double]"0.0001060937429429032"
$ws.Range("T5").Value = [double]"7.973251194925979e-05"
$ws.Range("U5").Value = [double]"0.0002137173141818494"
$ws.Range("V5").Value = [double]"1.540154698886909e-05"
$ws.Range("W5").Value = [double]"3.28467604049365e-06"
$ws.Range("X5").Value = [double]"0.0005501061677932739"
$ws.Range("Y5").Value = [double]"0.003062946023419499"

$ws.Range("A6").Value = [double]"0.0001661850401433185"
$ws.Range("B6").Value = [double]"0.0009311214089393616"
$ws.Range("C6").Value = [double]"1.177775210692289e-08"
$ws.Range("D6").Value = [double]"8.150353636438012e-08"
$ws.Range("E6").Value = [double]"0.005228939466178417"
$ws.Range("F6").Value = [double]"6.134509078492556e-08"
$ws.Range("G6").Value = [double]"0.0002242642513010651"
$ws.Range("H6").Value = [double]"1.149991734350726e-17"
$ws.Range("I6").Value = [double]"8.601014997111633e-06"
$ws.Range("J6").Value = [double]"0.0001464478846173733"
$ws.Range("K6").Value = [double]"7.411131809931248e-05"
$ws.Range("L6").Value = [double]"1.303019899518852e-17"
$ws.Range("M6").Value = [double]"0.006362977437674999"
$ws.Range("N6").Value = [double]"6.923304858083054e-18"
$ws.Range("O6").Value = [double]"5.276109277474461e-06"
$ws.Range("P6").Value = [double]"3.476714055539176e-11"
$ws.Range("Q6").Value = [double]"8.808006555227621e-07"
$ws.Range("R6").Value = [double]"5.007824256608728e-06"
$ws.Range("S6").Value = [double]"0.0002270957047585398"
$ws.Range("T6").Value = [double]"2.061329951175139e-08"
$ws.Range("U6").Value = [double]"7.879411495269822e-18"
$ws.Range("V6").Value = [double]"3.971583240058862e-08"
$ws.Range("W6").Value = [double]"3.122039800018683e-07"
$ws.Range("X6").Value = [double]"2.542059618360071e-11"
$ws.Range("Y6").Value = [double]"0.001226022955961525"

$ws.Range("A7").Value = [double]"0.0003369326295796782"
$ws.Range("B7").Value = [double]"0.0007517669582739472"
$ws.Range("C7").Value = [double]"0.0001107530333683826"
$ws.Range("D7").Value = [double]"1.853537838769626e-07"
$ws.Range("E7").Value = [double]"4.447271422918115e-19"
$ws.Range("F7").Value = [double]"1.81109285790626e-07"
$ws.Range("G7").Value = [double]"0.0002258610766148195"
$ws.Range("H7").Value = [double]"1.445805963973184e-16"
$ws.Range("I7").Value = [double]"0.000140322299557738"
$ws.Range("J7").Value = [double]"0.001625792006962001"
$ws.Range("K7").Value = [double]"3.267081183366827e-08"
$ws.Range("L7").Value = [double]"4.314128432270081e-07"
$ws.Range("M7").Value = [double]"0.009886804036796093"
$ws.Range("N7").Value = [double]"1.036427131980198e-16"
$ws.Range("O7").Value = [double]"4.240775524522178e-05"
$ws.Range("P7").Value = [double]"3.262393022396282e-07"
$ws.Range("Q7").Value = [double]"1.620000148250256e-05"
$ws.Range("R7").Value = [double]"6.430234122944967e-08"
$ws.Range("S7").Value = [double]"6.90019951434806e-05"
$ws.Range("T7").Value = [double]"9.146166121354327e-05"
$ws.Range("U7").Value = [double]"0.0002531820500735193"
$ws.Range("V7").Value = [double]"7.109832222340629e-05"
$ws.Range("W7").Value = [double]"5.794519438495627e-06"
$ws.Range("X7").Value = [double]"0.0006738266092725098"
$ws.Range("Y7").Value = [double]"0.002702496480196714"

$ws.Range("A8").Value = [double]"0.0001596811489434913"
$ws.Range("B8").Value = [double]"0.0007996290805749595"
$ws.Range("C8").Value = [double]"2.362860663287494e-15"
$ws.Range("D8").Value = [double]"1.256528776139021e-06"
$ws.Range("E8").Value = [double]"0.0003429708594921976"
$ws.Range("F8").Value = [double]"1.098644816011074e-06"
$ws.Range("G8").Value = [double]"5.37697305844631e-05"
$ws.Range("H8").Value = [double]"3.238285728457413e-07"
$ws.Range("I8").Value = [double]"9.783330464205165e-17"
$ws.Range("J8").Value = [double]"0.0003655441978480667"
$ws.Range("K8").Value = [double]"1.69156464835396e-06"
$ws.Range("L8").Value = [double]"7.461642326234141e-07"
$ws.Range("M8").Value = [double]"0.01067871227860451"
$ws.Range("N8").Value = [double]"5.583078705483812e-17"
$ws.Range("O8").Value = [double]"6.104653584770858e-05"
$ws.Range("P8").Value = [double]"6.301378107309574e-06"
$ws.Range("Q8").Value = [double]"1.820612283154333e-07"
$ws.Range("R8").Value = [double]"6.159193435451016e-05"
$ws.Range("S8").Value = [double]"0.0003235218173358589"
$ws.Range("T8").Value = [double]"0.0002318381448276341"
$ws.Range("U8").Value = [double]"1.43867546285037e-05"
$ws.Range("V8").Value = [double]"0.001272225752472878"
$ws.Range("W8").Value = [double]"2.197476169385482e-05"
$ws.Range("X8").Value = [double]"0.0007041794015094638"
$ws.Range("Y8").Value = [double]"0.002106880536302924"

$ws.Range("A9").Value = [double]"0.0002469881437718868"
$ws.Range("B9").Value = [double]"0.000550118216779083"
$ws.Range("C9").Value = [double]"1.912545667437371e-05"
$ws.Range("D9").Value = [double]"1.698408937045315e-06"
$ws.Range("E9").Value = [double]"1.843013164124697e-17"
$ws.Range("F9").Value = [double]"4.558618513783586e-07"
$ws.Range("G9").Value = [double]"0.0002102444559568539"
$ws.Range("H9").Value = [double]"1.928758540985821e-17"
$ws.Range("I9").Value = [double]"0.0001744547043927014"
$ws.Range("J9").Value = [double]"4.910009374725632e-05"
$ws.Range("K9").Value = [double]"1.374942826259939e-06"
$ws.Range("L9").Value = [double]"8.712789281162259e-07"
$ws.Range("M9").Value = [double]"0.009681799449026585"
$ws.Range("N9").Value = [double]"4.515666451876399e-10"
$ws.Range("O9").Value = [double]"0.0009832425275817513"
$ws.Range("P9").Value = [double]"1.136164442527843e-07"
$ws.Range("Q9").Value = [double]"0.0001140882959589362"
$ws.Range("R9").Value = [double]"3.828645500743336e-17"
$ws.Range("S9").Value = [double]"6.878614658489823e-05"
$ws.Range("T9").Value = [double]"9.424791642231867e-06"
$ws.Range("U9").Value = [double]"0.0003191560972481966"
$ws.Range("V9").Value = [double]"0.0002470989711582661"
$ws.Range("W9").Value = [double]"7.276165888470132e-06"
$ws.Range("X9").Value = [double]"1.400855489919195e-05"
$ws.Range("Y9").Value = [double]"0.00296163116581738"

$ws.Range("A10").Value = [double]"0.0004361461324151605"
$ws.Range("B10").Value = [double]"0.0008663510670885444"
$ws.Range("C10").Value = [double]"4.186919341009343e-06"
$ws.Range("D10").Value = [double]"2.181029709863013e-14"
$ws.Range("E10").Value = [double]"3.28313014380421e-17"
$ws.Range("F10").Value = [double]"0.001906165038235486"
$ws.Range("G10").Value = [double]"0.000521318637765944"
$ws.Range("H10").Value = [double]"3.800306954701759e-17"
$ws.Range("I10").Value = [double]"0.0002027807786362246"
$ws.Range("J10").Value = [double]"0.0003878700081259012"
$ws.Range("K10").Value = [double]"3.023836541160563e-07"
$ws.Range("L10").Value = [double]"2.357251167741968e-17"
$ws.Range("M10").Value = [double]"0.01354031451046467"
$ws.Range("N10").Value = [double]"4.912015529532612e-17"
$ws.Range("O10").Value = [double]"4.045211881020805e-07"
$ws.Range("P10").Value = [double]"6.320616100693243e-17"
$ws.Range("Q10").Value = [double]"0.000174307293491438"
$ws.Range("R10").Value = [double]"5.480202958096925e-07"
$ws.Range("S10").Value = [double]"0.0002562867302913219"
$ws.Range("T10").Value = [double]"0.0002257677115267143"
$ws.Range("U10").Value = [double]"0.0002562968293204904"
$ws.Range("V10").Value = [double]"7.58857058826834e-05"
$ws.Range("W10").Value = [double]"4.04053116653813e-06"
$ws.Range("X10").Value = [double]"0.0003391627687960863"
$ws.Range("Y10").Value = [double]"0.003512072376906872"

$ws.Range("A11").Value = [double]"0.0002580983273219317"
$ws.Range("B11").Value = [double]"0.0007403444615192711"
$ws.Range("C11").Value = [double]"5.471041731652804e-06"
$ws.Range("D11").Value = [double]"2.182753178203711e-06"
$ws.Range("E11").Value = [double]"0.002939575584605336"
$ws.Range("F11").Value = [double]"3.47649847753928e-06"
$ws.Range("G11").Value = [double]"0.0001829219982028008"
$ws.Range("H11").Value = [double]"2.422149442497163e-18"
$ws.Range("I11").Value = [double]"0.0001211178459925577"
$ws.Range("J11").Value = [double]"9.777095328900032e-06"
$ws.Range("K11").Value = [double]"3.242481398046948e-06"
$ws.Range("L11").Value = [double]"5.129232594214272e-09"
$ws.Range("M11").Value = [double]"0.009115278720855713"
$ws.Range("N11").Value = [double]"9.80828915026385e-18"
$ws.Range("O11").Value = [double]"3.097480374100769e-18"
$ws.Range("P11").Value = [double]"9.952897380571812e-05"
$ws.Range("Q11").Value = [double]"3.237030614400283e-05"
$ws.Range("R11").Value = [double]"2.22148855755222e-06"
$ws.Range("S11").Value = [double]"0.0003115073195658624"
$ws.Range("T11").Value = [double]"1.167744630947709e-05"
$ws.Range("U11").Value = [double]"0.000315976794809103"
$ws.Range("V11").Value = [double]"4.540929876384325e-05"
$ws.Range("W11").Value = [double]"0.0001954954495886341"
$ws.Range("X11").Value = [double]"0.002216053893789649"
$ws.Range("Y11").Value = [double]"0.0008682330953888595"

Write-Output "Updated 250 cells in rows 2-11"